# Remove "Marky Marc", and add flight details to reports.
#
# The guest-pair row that used to read  A30=Marky Marc / B30=Trevor Anthony
# loses "Marky Marc": Trevor Anthony moves up into Guest1 (A30) and a brand
# new guest, "Manuel Daquilema", is placed into Guest2 (B30). The room is no
# longer confirmed, so the Confirmed flag (E30) flips from 1 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A30").Value = "Trevor Anthony"
$ws.Range("B30").Value = "Manuel Daquilema"
$ws.Range("E30").Value = 0

# Reposition the view the way the author left it before saving: scrolled
# down so row 22 is at the top, with B30 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("B30").Select()
